$d = $word.ActiveDocument

# 1) Remove the first four paragraphs:
#    "Note: This is not finalized Syllabus Batch One -> ghp_t5CV..."
#    (empty paragraph)
#    "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"
#    (empty paragraph)
$p1 = $d.Paragraphs.Item(1)
$p5 = $d.Paragraphs.Item(5)
$introRange = $d.Range($p1.Range.Start, $p5.Range.Start)
$introRange.Delete()

# 2) Remove the stray lastRenderedPageBreak marker that precedes "Project 5".
#    Re-serializing the run's own OOXML and re-inserting it in place drops the
#    transient lastRenderedPageBreak marker (it is not part of a sub-range's
#    canonical XML) while preserving all real run/paragraph formatting.
$find = $d.Content.Find
$find.ClearFormatting()
$found = $find.Execute("Project 5", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $target = $d.Content.Duplicate
    $target.Find.Execute("Project 5")
    $cleanXml = $target.WordOpenXML
    $target.Delete()
    $target.InsertXML($cleanXml)
}
